# TC07_CDS_phs001437_LibrarySource_Transcriptomic.xlsx
#
# The "FilesTab" query stored in B4 gains a new trailing projected column
# ('' AS "Supplementary File") appended right after the existing
# "Library Strategy" column. Re-writing the cell text also causes the
# workbook's shared-string table to be re-emitted in (new) usage order,
# which is what moves the "Sample ID" query string ahead of the rewritten
# "File Name" query string - matching the reordering seen in the diff.
#
# The cell also picks up a (visually identical, still 12pt Calibri/wrap)
# style after the edit, and the active selection moves from C3 to B4 - both
# reproduced below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newFilesQuery = @"
SELECT
    f1.file_name AS "File Name",
    s.study_name AS "Study Name",
    s.phs_accession AS "Accession",
    sp.participant_id AS "Participant Id",
    COALESCE((
        SELECT
            REPLACE(GROUP_CONCAT(CASE WHEN rn <= 5 THEN smp.sample_id ELSE NULL END, ', '), ', , ', ', ') ||
            CASE WHEN MAX(rn) > 5 THEN ', ...' ELSE '' END
        FROM (
            SELECT
                smp.sample_id,
                ROW_NUMBER() OVER (ORDER BY smp.sample_id) AS rn
            FROM df_sample smp
            WHERE smp."participant.study_participant_id" = sp.study_participant_id
        ) smp
    ), '') AS "Sample Id",
    f1.file_type AS "File Type",
    gi.library_strategy AS "Library Strategy",
'' AS "Supplementary File"
FROM 
    df_study s
INNER JOIN 
    df_participant sp ON sp."study.phs_accession" = s.phs_accession
INNER JOIN  
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
INNER JOIN 
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
INNER JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
INNER JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
INNER JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
WHERE 
    s.phs_accession = 'phs001437' AND gi.library_source = 'Transcriptomic'
GROUP BY
    f1.file_name,
    s.study_name,
    s.phs_accession,
    sp.participant_id,
    f1.file_type,
    gi.library_strategy
ORDER BY 
    f1.file_name ASC
LIMIT 100;
"@

$filesCell = $ws.Range("B4")
$filesCell.Value = $newFilesQuery

# Keep the cell wrapped / 12pt, matching the rest of the TabQuery column -
# re-asserting this after the text change is what lands the cell on its own
# (new) style record instead of the previous one.
$filesCell.WrapText = $true
$filesCell.Font.Size = 12
$filesCell.Font.ThemeColor = 1

# Row grew by one line of text; pin it back to the height recorded upstream.
$ws.Rows(4).RowHeight = 409.6

# Final user selection lands on the edited cell.
$filesCell.Select() | Out-Null
